$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.051542401313782
$ws.Range("B1").Value = 2.343771696090698
$ws.Range("C1").Value = 6.691954135894775
$ws.Range("D1").Value = 2.281288385391235
$ws.Range("E1").Value = 1.303734183311462
